# ---------------------------------------------------------------------------
# Applies the "WRI China / Hong Kong EPS v2.0.0" update to
# "Frac of Bldgs Owned by Entity.xlsx":
#   1. Insert a new "Output by Industry" sheet (between Commercial & FoBObE)
#      holding the GDP output-share table used to split commercial buildings
#      owned by "industry" into energy vs. non-energy suppliers.
#   2. Add workbook-level defined names pointing at that table.
#   3. Rework the FoBObE sheet: rename the title & some rows, and add new
#      rows for the disaggregated energy-supplier categories.
#   4. Add explanatory "Notes" text to the About sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "Output by Industry" sheet
# ---------------------------------------------------------------------------
$fobobe = $wb.Worksheets.Item("FoBObE")
$industry = $wb.Worksheets.Add($fobobe)
$industry.Name = "Output by Industry"

$industry.Columns.Item(1).ColumnWidth = 12
$industry.Columns.Item(2).ColumnWidth = 36

$industry.Range("A1").Value = "For bibliographic source and methods, see file output_shares_by_industry.xslx"
$industry.Range("A2").Value = "in the InputData folder."

$industry.Range("A4").Value = "Share"
$industry.Range("B4").Value = "Industry Category"
$industry.Range("A4:B4").Font.Bold = $true
$industry.Range("A4:B4").Interior.ColorIndex = 15

$industry.Range("A5").Value = 0.93219819361870848
$industry.Range("B5").Value = "non-energy industries"

$industry.Range("A6").Value = [double]"1.5490855293616566E-2"
$industry.Range("B6").Value = "electricity suppliers"

$industry.Range("A7").Value = [double]"2.1866536828369144E-3"
$industry.Range("B7").Value = "coal suppliers"

$industry.Range("A8").Value = [double]"4.9443136381930888E-2"
$industry.Range("B8").Value = "natural gas and petroleum suppliers"

$industry.Range("A9").Value = [double]"6.8116102290716575E-4"
$industry.Range("B9").Value = "biomass and biofuel suppliers"

$industry.Range("A10").Value = 0
$industry.Range("B10").Value = "other energy suppliers"

$industry.Range("A5:A9").NumberFormat = "0.0000"

# ---------------------------------------------------------------------------
# 2. Workbook-level defined names used by the FoBObE formulas below
# ---------------------------------------------------------------------------
$wb.Names.Add('outputfrac_nonenergy', "='Output by Industry'!`$A`$5")
$wb.Names.Add('outputfrac_elec',      "='Output by Industry'!`$A`$6")
$wb.Names.Add('outputfrac_coal',      "='Output by Industry'!`$A`$7")
$wb.Names.Add('outputfrac_ngps',      "='Output by Industry'!`$A`$8")
$wb.Names.Add('outputfrac_bio',       "='Output by Industry'!`$A`$9")
$wb.Names.Add('outputfrac_other',     "='Output by Industry'!`$A`$10")

# ---------------------------------------------------------------------------
# 3. FoBObE sheet rework
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FoBObE")

# Title + column headers
$ws.Range("A1").Value = "Ownership by Cash Flow Entity (dimensionless)"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30
$ws.Columns.Item(1).ColumnWidth = 35.28515625

# Rename existing category rows
$ws.Range("A2").Value = "government"
$ws.Range("A3").Value = "nonenergy industries"
$ws.Range("A4").Value = "labor and consumers"

# Row 3 commercial column now pulls its share from the new output table
$ws.Range("D3").Formula = '=Commercial!$C$21*outputfrac_nonenergy'

# New rows for the disaggregated energy-supplier categories
$ws.Range("A5").Value = "foreign entities"
$ws.Range("B5").Value = 0
$ws.Range("C5").Formula = "=B5"
$ws.Range("D5").Value = 0

$ws.Range("A6").Value = "electricity suppliers"
$ws.Range("B6").Value = 0
$ws.Range("C6").Formula = "=B6"
$ws.Range("D6").Formula = '=Commercial!$C$21*outputfrac_elec'

$ws.Range("A7").Value = "coal suppliers"
$ws.Range("B7").Value = 0
$ws.Range("C7").Formula = "=B7"
$ws.Range("D7").Formula = '=Commercial!$C$21*outputfrac_coal'

$ws.Range("A8").Value = "natural gas and petroleum suppliers"
$ws.Range("B8").Value = 0
$ws.Range("C8").Formula = "=B8"
$ws.Range("D8").Formula = '=Commercial!$C$21*outputfrac_ngps'

$ws.Range("A9").Value = "biomass and biofuel suppliers"
$ws.Range("B9").Value = 0
$ws.Range("C9").Formula = "=B9"
$ws.Range("D9").Formula = '=Commercial!$C$21*outputfrac_bio'

$ws.Range("A10").Value = "other energy suppliers"
$ws.Range("B10").Value = 0
$ws.Range("C10").Formula = "=B10"
$ws.Range("D10").Formula = '=Commercial!$C$21*outputfrac_other'

# ---------------------------------------------------------------------------
# 4. About sheet notes
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A18").Value = "Notes"
$about.Range("A18").Font.Bold = $true

$about.Range("A19").Value = "Fraction of buildings owned by cash flow entity is used in cash flow calculations."
$about.Range("A20").Value = 'For commercial buildings, we divde between "nonenergy industries" and'
$about.Range("A21").Value = "the various energy industries based on output shares of GDP.  (We assume"
$about.Range("A22").Value = 'that all residential buildings owned by industry are owned by "nonenergy'
$about.Range("A23").Value = 'industries" - i.e. rental property management communities, nursing homes,'
$about.Range("A24").Value = "etc. - not energy industries."

$about.Range("A26").Value = "In the output tab, we show more decimal places than the source data"
$about.Range("A27").Value = "provide in order to avoid rounding error in Vensim (each column must"
$about.Range("A28").Value = "add to 1)."

$about.Activate()
$about.Range("B25").Select()

$wb.Application.CalculateFull()
